$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 11750
$ws.Range("J3").Value = 11750
$ws.Range("L3").Value = 11750
$ws.Range("N3").Value = -11978
$ws.Range("H19").Value = 1686.875
$ws.Range("I19").Value = 1599.2
$ws.Range("K19").Value = 1599.2
$ws.Range("M19").Value = -1424.2
$ws.Range("H33").Value = 239.54546
$ws.Range("I33").Value = 245.95238
$ws.Range("J33").Value = 105
$ws.Range("K33").Value = 245.95238
$ws.Range("L33").Value = 105
$ws.Range("M33").Value = -16.95238000000001
$ws.Range("N33").Value = -563
$ws.Range("H64").Value = 4000
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 4000
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716
$ws.Range("H102").Value = 11750
$ws.Range("J102").Value = 11750
$ws.Range("L102").Value = 11750
$ws.Range("N102").Value = -18240
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23147.857
$ws.Range("I32").Value = 23147.857
$ws.Range("K32").Value = 23147.857
$ws.Range("M32").Value = -22860.857
$ws.Range("H45").Value = 3150
$ws.Range("I45").Value = 3150
$ws.Range("K45").Value = 3150
$ws.Range("M45").Value = -2773
$ws.Range("H61").Value = 2158.3333
$ws.Range("I61").Value = 2158.3333
$ws.Range("K61").Value = 2158.3333
$ws.Range("M61").Value = -1946.3333
$ws.Range("H74").Value = 2364.4375
$ws.Range("I74").Value = 1563.9231
$ws.Range("K74").Value = 1563.9231
$ws.Range("M74").Value = -689.9231
$ws.Range("H76").Value = 52414.285
$ws.Range("J76").Value = 52414.285
$ws.Range("L76").Value = 52414.285
$ws.Range("N76").Value = -53090.285
$ws.Range("H77").Value = 2364.4375
$ws.Range("I77").Value = 1563.9231
$ws.Range("K77").Value = 7819.6155
$ws.Range("M77").Value = -3451.6155
$ws.Range("H79").Value = 52414.285
$ws.Range("J79").Value = 52414.285
$ws.Range("L79").Value = 52414.285
$ws.Range("N79").Value = -54754.285
$ws.Range("H136").Value = 2158.3333
$ws.Range("I136").Value = 2158.3333
$ws.Range("K136").Value = 6474.999899999999
$ws.Range("M136").Value = -3924.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 698.75
$ws.Range("I22").Value = 599.3333
$ws.Range("K22").Value = 599.3333
$ws.Range("M22").Value = -426.3333
$ws.Range("H88").Value = 19633.666
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 19633.666
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 19633.666
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -20445.666
$ws.Range("H91").Value = 19633.666
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 19633.666
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 19633.666
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -22441.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 47787.25
$ws.Range("J28").Value = 47787.25
$ws.Range("L28").Value = 47787.25
$ws.Range("N28").Value = -48277.25
$ws.Range("H31").Value = 3773.5715
$ws.Range("I31").Value = 2491.6667
$ws.Range("K31").Value = 2491.6667
$ws.Range("M31").Value = -2196.6667
$ws.Range("H34").Value = 3773.5715
$ws.Range("I34").Value = 2491.6667
$ws.Range("K34").Value = 2491.6667
$ws.Range("M34").Value = -2289.6667
$ws.Range("H58").Value = 1824.875
$ws.Range("J58").Value = 2595
$ws.Range("L58").Value = 2595
$ws.Range("N58").Value = -3001
$ws.Range("H74").Value = 8000
$ws.Range("J74").Value = 8000
$ws.Range("L74").Value = 8000
$ws.Range("N74").Value = -9748
$ws.Range("H77").Value = 8000
$ws.Range("J77").Value = 8000
$ws.Range("L77").Value = 24000
$ws.Range("N77").Value = -32736
$ws.Range("H88").Value = 79000
$ws.Range("J88").Value = 79000
$ws.Range("L88").Value = 79000
$ws.Range("N88").Value = -79812
$ws.Range("H91").Value = 79000
$ws.Range("J91").Value = 79000
$ws.Range("L91").Value = 79000
$ws.Range("N91").Value = -81808
$ws.Range("H103").Value = 32850
$ws.Range("I103").Value = 32850
$ws.Range("K103").Value = 32850
$ws.Range("M103").Value = -31678
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 4111.75
$ws.Range("I132").Value = 4111.75
$ws.Range("K132").Value = 12335.25
$ws.Range("M132").Value = -9805.25
$ws.Range("H136").Value = 1824.875
$ws.Range("J136").Value = 2595
$ws.Range("L136").Value = 7785
$ws.Range("N136").Value = -12885

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 381.25
$ws.Range("I35").Value = 341.66666
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 1024.99998
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = -736.99998
$ws.Range("N35").Value = -2076
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H97").Value = 750
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5016.5
$ws.Range("I33").Value = 5016.5
$ws.Range("K33").Value = 5016.5
$ws.Range("M33").Value = -4764.5
$ws.Range("H101").Value = 48332.668
$ws.Range("J101").Value = 48332.668
$ws.Range("L101").Value = 48332.668
$ws.Range("N101").Value = -54822.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 750
$ws.Range("I46").Value = 750
$ws.Range("K46").Value = 750
$ws.Range("M46").Value = -562
$ws.Range("H95").Value = 29000
$ws.Range("J95").Value = 29000
$ws.Range("L95").Value = 29000
$ws.Range("N95").Value = -34492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 27241.5
$ws.Range("J69").Value = 27241.5
$ws.Range("L69").Value = 27241.5
$ws.Range("N69").Value = -28739.5
$ws.Range("H72").Value = 27241.5
$ws.Range("J72").Value = 27241.5
$ws.Range("L72").Value = 81724.5
$ws.Range("N72").Value = -89212.5
